# Generate Report for Handback
#
# The localization-status workbook previously recorded each localized
# file as "Ready for handoff". This run marks the .md source file in
# every language sheet as handed back (in sync with en-US) and records
# the handback report: the target (handoff) file, the handback file
# (the same translated artifact being returned) and the datetime the
# handback happened.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: just flips the status label for both languages ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusText
$overview.Range("C2").Value = $statusText

# Hyperlink font look-alike of the workbook's existing custom "HyperLink"
# cell style (underline, Calibri 11, RGB 6495ED stored BGR-packed for
# the Font.Color COM property).
$hyperlinkFontColor = 15570276

function Set-HandbackRow {
    param($ws, $handbackDateTime)

    # B2: status column -> handed back
    $ws.Range("B2").Value = $statusText

    # Find the existing hyperlinks on the source (.md, column A) and the
    # latest handoff file (.xlf, column C) for row 2 so the new "target"
    # and "handback" columns can point at the very same files.
    $mdLink = $null
    $xlfLink = $null
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Row -eq 2 -and $hl.Range.Column -eq 1) { $mdLink = $hl }
        if ($hl.Range.Row -eq 2 -and $hl.Range.Column -eq 3) { $xlfLink = $hl }
    }

    # E2: Latest Target File -> same source file as column A
    $ws.Hyperlinks.Add($ws.Range("E2"), $mdLink.Address, "", "", $mdLink.TextToDisplay)
    $ws.Range("E2").Font.Underline = $true
    $ws.Range("E2").Font.Color = $hyperlinkFontColor
    $ws.Range("E2").Font.Name = "Calibri"
    $ws.Range("E2").Font.Size = 11

    # F2: Latest Handback File -> same translated artifact as column C
    $ws.Hyperlinks.Add($ws.Range("F2"), $xlfLink.Address, "", "", $xlfLink.TextToDisplay)
    $ws.Range("F2").Font.Underline = $true
    $ws.Range("F2").Font.Color = $hyperlinkFontColor
    $ws.Range("F2").Font.Name = "Calibri"
    $ws.Range("F2").Font.Size = 11

    # G2: Latest Handback DateTime -> when the handback report was produced
    $ws.Range("G2").Value = $handbackDateTime
}

$zhcn = $wb.Worksheets.Item("zh-cn")
Set-HandbackRow $zhcn "2016-02-22 12:26:20"

$dede = $wb.Worksheets.Item("de-de")
Set-HandbackRow $dede "2016-02-22 12:26:44"
